$wb = $excel.ActiveWorkbook

# --- gc_fields_display: rename "best rolling" fields to "critical" ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A82").Value = "__CalcBestRollingWeightedMeanHeartRate"
$ws1.Range("B82").Value = "Critical Heart Rate"

$ws1.Range("A83").Value = "__CalcBestRollingWeightedMeanPower"
$ws1.Range("B83").Value = "Critical Power"

$ws1.Range("A84").Value = "__CalcBestRollingWeightedMeanSpeed"
$ws1.Range("B84").Value = "Critical Speed"

$ws1.Range("A85").Value = "__CalcBestRollingWeightedMeanPace"
$ws1.Range("B85").Value = "Critical Pace"

# --- gc_fields_uom: add matching unit-of-measure rows + fix a statute unit ---
$ws2 = $wb.Worksheets.Item(2)

# BeginPowerTwentyMinutesDistance statute unit should be mile, not kilometer
$ws2.Range("D2").Value = "mile"

$ws2.Range("A46").Value = "__CalcBestRollingWeightedMeanHeartRate"
$ws2.Range("B46").Value = "all"
$ws2.Range("C46").Value = "bpm"
$ws2.Range("D46").Value = "bpm"

$ws2.Range("A47").Value = "__CalcBestRollingWeightedMeanPower"
$ws2.Range("B47").Value = "all"
$ws2.Range("C47").Value = "watt"
$ws2.Range("D47").Value = "watt"

$ws2.Range("A48").Value = "__CalcBestRollingWeightedMeanSpeed"
$ws2.Range("B48").Value = "all"
$ws2.Range("C48").Value = "kph"
$ws2.Range("D48").Value = "mph"

$ws2.Range("A49").Value = "__CalcBestRollingWeightedMeanPace"
$ws2.Range("B49").Value = "all"
$ws2.Range("C49").Value = "minperkm"
$ws2.Range("D49").Value = "minpermile"

# --- view state: gc_fields_uom becomes the active/selected sheet ---
$ws1.Range("A82:A85").Select()
$ws2.Range("C50").Select()
$ws2.Activate()
